# Add invoicing periods bounds to example 08
#
# Summary of the edit (per the commit "Add invoicing periods bounds to
# example 08"):
#   - "experts" sheet: selection moves from C1 to A3
#   - "invoicing periods bounds" sheet:
#       * C2:C4 lower bound changed from 100 to 0
#       * D7:D11 upper bound changed to a flat 200 (was 202..206)
#       * nine new data rows (12..20) appended, one per remaining expert,
#         each with Lower=0 / Upper=200 and the usual COUNTIF checks
#       * becomes the active / selected sheet, selection moves to D4
#   - workbook re-activates on "invoicing periods bounds" (tab 9)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "experts" sheet: just a selection change (C1 -> A3)
# ---------------------------------------------------------------------
$wsExperts = $wb.Worksheets.Item("experts")
$wsExperts.Range("A3").Select()

# ---------------------------------------------------------------------
# "invoicing periods bounds" sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("invoicing periods bounds")

# Lower bound 100 -> 0 for rows 2-4
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 0

# Upper bound normalised to 200 for rows 7-11 (was 202,203,204,205,206)
$ws.Range("D7").Value = 200
$ws.Range("D8").Value = 200
$ws.Range("D9").Value = 200
$ws.Range("D10").Value = 200
$ws.Range("D11").Value = 200

# New rows 12-20: one per expert, mirroring rows 7-11's layout/styles
$names = @("DEV.Barłomiej", "DEV.Cezary", "DEV.Dariusz", "DEV.Eugenius", "DEV.Franciszek", "DEV.Gustaw", "DEV.Hubert", "DEV.Ignacy", "DEV.Jarosław")

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = 12 + $i
    $name = $names[$i]

    $ws.Range("A$row").Value = $name
    $ws.Range("A$row").Style = $ws.Range("A11").Style

    $ws.Range("B$row").Value = "25.Jan"
    $ws.Range("B$row").Style = $ws.Range("B11").Style

    $ws.Range("C$row").Value = 0
    $ws.Range("C$row").Style = $ws.Range("C11").Style

    $ws.Range("D$row").Value = 200
    $ws.Range("D$row").Style = $ws.Range("D11").Style

    $ws.Range("E$row").Formula = "=COUNTIF(experts!`$A`$2:`$A`$985, A$row) > 0"
    $ws.Range("E$row").Style = $ws.Range("E11").Style

    $ws.Range("F$row").Formula = "=COUNTIF('invoicing periods'!`$A`$2:`$A`$998, B$row) > 0"
    $ws.Range("F$row").Style = $ws.Range("F11").Style
}

# Make this the active sheet/tab and move the selection to D4
$ws.Activate()
$ws.Range("D4").Select()
